# Generate Report for Handback
#
# For each localized-language sheet (zh-cn, de-de):
#   - Status (col C) rows 2-3: "Ready for handoff" -> "Handed back: in sync with en-US"
#   - Latest Target File (col F) rows 2-3: filled with the source .md file name,
#     hyperlinked to the same target as column A's hyperlink on that row.
#   - Latest Handback File (col G) rows 2-3: filled with the .xlf file name,
#     hyperlinked to the same target as column D's hyperlink on that row.
#   - Latest Handback DateTime (col H) rows 2-3: "0001-01-01 00:00:00" -> actual
#     handback timestamp (same value for both rows on a given sheet, distinct
#     per sheet/locale).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$sheets = @(
    @{
        Name = "zh-cn"
        HandbackDateTime = "2016-03-17 10:26:25"
        MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/f3ecc755f260fc6836a636eebb207ab60af5593d/e2e"
        XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/265b2c67abacd95070db083ae99917fa40db4a21/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
    },
    @{
        Name = "de-de"
        HandbackDateTime = "2016-03-17 10:26:30"
        MdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/f3ecc755f260fc6836a636eebb207ab60af5593d/e2e"
        XlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/811854fef82c7b39bc1827d1cd544e7b53af3d40/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"
    }
)

$rows = @(
    @{ Row = 2; Md = "db1c8fb5-6c61-4671-b979-b1de64768167.md"; Xlf = "db1c8fb5-6c61-4671-b979-b1de64768167.177ad54f80913b18d873e0b194c4680698bfb610" },
    @{ Row = 3; Md = "ee64cd91-92b5-4a89-82c4-17dde9f1fac5.md"; Xlf = "ee64cd91-92b5-4a89-82c4-17dde9f1fac5.2c7c27e890347422a65c36252ed763141bd3e984" }
)

foreach ($sheetInfo in $sheets) {
    $ws = $wb.Worksheets.Item($sheetInfo.Name)

    foreach ($r in $rows) {
        $rowNum = $r.Row

        # Status -> Handed back
        $ws.Cells.Item($rowNum, 3).Value = $statusText

        # Latest Target File (F) - same file/link as column A (source .md)
        $mdUrl = $sheetInfo.MdAddress + "/" + $r.Md
        $ws.Hyperlinks.Add($ws.Cells.Item($rowNum, 6), $mdUrl, "", "", $r.Md)

        # Latest Handback File (G) - same file/link as column D (.xlf target file)
        $xlfFile = $r.Xlf + "." + $sheetInfo.Name + ".xlf"
        $xlfUrl = $sheetInfo.XlfAddress + "/" + $xlfFile
        $ws.Hyperlinks.Add($ws.Cells.Item($rowNum, 7), $xlfUrl, "", "", $xlfFile)

        # Latest Handback DateTime (H)
        $ws.Cells.Item($rowNum, 8).Value = $sheetInfo.HandbackDateTime
    }
}
